# Auto-generated edit script updating cryptos list prices/volumes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.199.85"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.480.73"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.28"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.35"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.44%  "

$ws.Range("E9").Value = "  +2.79%  "

$ws.Range("E10").Value = "  +0.96%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.96"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.33%  "

$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.47"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.095.77"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("E16").Value = "  +0.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.467.73"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.61"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.95"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "350.81"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.03%  "

$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("E22").Value = "  -0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.86"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.23"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("E25").Value = "  +2.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.19"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.607.49"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0910"
$ws.Range("D29").ClearFormats()

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "504.36"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.01%  "

$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("E32").Value = "  +0.57%  "

$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("E35").Value = "  +1.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.87"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.70"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.20"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.03%  "

$ws.Range("E39").Value = "  -1.02%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.70"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.88%  "

$ws.Range("E42").Value = "  +0.57%  "

$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.39"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "143.09"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0263"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.49"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.90%  "

$ws.Range("E48").Value = "  +0.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0740"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.39%  "

$ws.Range("E51").Value = "  +0.55%  "

